# Append the 19 new daily rows (43-61) of COVID tracking data to the 'aglomerados' sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (dates): reuse the m/d/yy date format already applied to rows 2-42.
$dateSerials = @(43983, 43984, 43985, 43986, 43987, 43988, 43989, 43990, 43991, 43992, 43993, 43994, 43995, 43996, 43997, 43998, 43999, 44000, 44001)
$startRow = 43

for ($i = 0; $i -lt $dateSerials.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $dateSerials[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "m/d/yy"
}

$data = New-Object 'object[,]' 19,12
$data[0,0] = 2711
$data[0,1] = 2636
$data[0,2] = 41
$data[0,3] = 164
$data[0,4] = 749
$data[0,5] = 67
$data[0,6] = 176
$data[0,7] = 0.22800000000000001
$data[0,8] = 761
$data[0,9] = [double]"7.3999999999999996E-2"
$data[0,10] = 83
$data[0,11] = 81
$data[1,0] = 2862
$data[1,1] = 2817
$data[1,2] = 46
$data[1,3] = 188
$data[1,4] = 799
$data[1,5] = 75
$data[1,6] = 162
$data[1,7] = 0.30199999999999999
$data[1,8] = 772
$data[1,9] = [double]"7.2999999999999995E-2"
$data[1,10] = 86
$data[1,11] = 102
$data[2,0] = 2916
$data[2,1] = 3029
$data[2,2] = 46
$data[2,3] = 177
$data[2,4] = 929
$data[2,5] = 79
$data[2,6] = 161
$data[2,7] = 0.29399999999999998
$data[2,8] = 744
$data[2,9] = [double]"7.0999999999999994E-2"
$data[2,10] = 77
$data[2,11] = 100
$data[3,0] = 3153
$data[3,1] = 3388
$data[3,2] = 50
$data[3,3] = 221
$data[3,4] = 1065
$data[3,5] = 84
$data[3,6] = 150
$data[3,7] = 0.34200000000000003
$data[3,8] = 733
$data[3,9] = [double]"8.5000000000000006E-2"
$data[3,10] = 97
$data[3,11] = 124
$data[4,0] = 3353
$data[4,1] = 3568
$data[4,2] = 53
$data[4,3] = 232
$data[4,4] = 1145
$data[4,5] = 91
$data[4,6] = 258
$data[4,7] = 0.376
$data[4,8] = 816
$data[4,9] = 0.113
$data[4,10] = 118
$data[4,11] = 114
$data[5,0] = 3466
$data[5,1] = 3788
$data[5,2] = 53
$data[5,3] = 238
$data[5,4] = 1203
$data[5,5] = 101
$data[5,6] = 258
$data[5,7] = 0.36299999999999999
$data[5,8] = 816
$data[5,9] = 0.11600000000000001
$data[5,10] = 121
$data[5,11] = 117
$data[6,0] = 3517
$data[6,1] = 4033
$data[6,2] = 56
$data[6,3] = 240
$data[6,4] = 1243
$data[6,5] = 113
$data[6,6] = 228
$data[6,7] = 0.45600000000000002
$data[6,8] = 816
$data[6,9] = 0.11700000000000001
$data[6,10] = 118
$data[6,11] = 122
$data[7,0] = 3633
$data[7,1] = 4243
$data[7,2] = 59
$data[7,3] = 227
$data[7,4] = 1454
$data[7,5] = 126
$data[7,6] = 228
$data[7,7] = 0.47399999999999998
$data[7,8] = 816
$data[7,9] = 0.115
$data[7,10] = 102
$data[7,11] = 125
$data[8,0] = 3848
$data[8,1] = 4504
$data[8,2] = 65
$data[8,3] = 246
$data[8,4] = 1597
$data[8,5] = 140
$data[8,6] = 228
$data[8,7] = 0.69299999999999995
$data[8,8] = 816
$data[8,9] = 0.16300000000000001
$data[8,10] = 117
$data[8,11] = 129
$data[9,0] = 4170
$data[9,1] = 4762
$data[9,2] = 66
$data[9,3] = 251
$data[9,4] = 1697
$data[9,5] = 153
$data[9,6] = 228
$data[9,7] = 0.64900000000000002
$data[9,8] = 816
$data[9,9] = 0.16200000000000001
$data[9,10] = 118
$data[9,11] = 133
$data[10,0] = 4291
$data[10,1] = 5086
$data[10,2] = 67
$data[10,3] = 244
$data[10,4] = 1792
$data[10,5] = 163
$data[10,6] = 233
$data[10,7] = 0.747
$data[10,8] = 816
$data[10,9] = 0.18099999999999999
$data[10,10] = 108
$data[10,11] = 136
$data[11,0] = 4497
$data[11,1] = 5390
$data[11,2] = 68
$data[11,3] = 237
$data[11,4] = 1912
$data[11,5] = 181
$data[11,6] = 233
$data[11,7] = 0.69499999999999995
$data[11,8] = 816
$data[11,9] = 0.17899999999999999
$data[11,10] = 107
$data[11,11] = 130
$data[12,0] = 4625
$data[12,1] = 5739
$data[12,2] = 68
$data[12,3] = 245
$data[12,4] = 2042
$data[12,5] = 199
$data[12,6] = 233
$data[12,7] = 0.70799999999999996
$data[12,8] = 816
$data[12,9] = 0.187
$data[12,10] = 118
$data[12,11] = 127
$data[13,0] = 4781
$data[13,1] = 6108
$data[13,2] = 70
$data[13,3] = 270
$data[13,4] = 2117
$data[13,5] = 213
$data[13,6] = 233
$data[13,7] = 0.746
$data[13,8] = 816
$data[13,9] = 0.217
$data[13,10] = 140
$data[13,11] = 130
$data[14,0] = 5102
$data[14,1] = 6390
$data[14,2] = 71
$data[14,3] = 296
$data[14,4] = 2386
$data[14,5] = 223
$data[14,6] = 233
$data[14,7] = 0.75900000000000001
$data[14,8] = 816
$data[14,9] = 0.23300000000000001
$data[14,10] = 145
$data[14,11] = 151
$data[15,0] = 5391
$data[15,1] = 6877
$data[15,2] = 72
$data[15,3] = 310
$data[15,4] = 2401
$data[15,5] = 250
$data[15,6] = 243
$data[15,7] = 0.74
$data[15,8] = 816
$data[15,9] = 0.22500000000000001
$data[15,10] = 144
$data[15,11] = 166
$data[16,0] = 5712
$data[16,1] = 7361
$data[16,2] = 77
$data[16,3] = 346
$data[16,4] = 2618
$data[16,5] = 272
$data[16,6] = 248
$data[16,7] = 0.76600000000000001
$data[16,8] = 816
$data[16,9] = 0.21
$data[16,10] = 166
$data[16,11] = 180
$data[17,0] = 6174
$data[17,1] = 8166
$data[17,2] = 80
$data[17,3] = 367
$data[17,4] = 2766
$data[17,5] = 295
$data[17,6] = 255
$data[17,7] = 0.8
$data[17,8] = 816
$data[17,9] = 0.28000000000000003
$data[17,10] = 185
$data[17,11] = 182
$data[18,0] = 6461
$data[18,1] = 8765
$data[18,2] = 84
$data[18,3] = 399
$data[18,4] = 3010
$data[18,5] = 322
$data[18,6] = 256
$data[18,7] = 0.82399999999999995
$data[18,8] = 816
$data[18,9] = 0.218
$data[18,10] = 219
$data[18,11] = 180

$ws.Range("B43:M61").Value2 = $data

# Row 48 picked up an explicit row height in the source workbook (matches the authored diff).
$ws.Rows.Item(48).RowHeight = 15.75

# Restore the view: scrolled so row 57 is the top-visible row, with N61 the active selection.
$ws.Range("N61").Select()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 2

Write-Output "Added rows $startRow-$($startRow + $dateSerials.Length - 1)"
